$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 111964457
$ws.Range("B2").Value = 56398
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = 'Tretåig hackspett'
$ws.Range("G2").Value = 'Picoides tridactylus'
$ws.Range("H2").Value = '(Linnaeus, 1758)'
$ws.Range("Q2").Value = 734949.4564622594
$ws.Range("R2").Value = 7088268.525185317
$ws.Range("Z2").Value = '16:01'
$ws.Range("AB2").Value = '16:01'

# Row 3
$ws.Range("A3").Value = 111965439
$ws.Range("B3").Value = 56398
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = 'Tretåig hackspett'
$ws.Range("G3").Value = 'Picoides tridactylus'
$ws.Range("H3").Value = '(Linnaeus, 1758)'
$ws.Range("Q3").Value = 734926.7697699566
$ws.Range("R3").Value = 7088234.05367971
$ws.Range("Z3").Value = '16:40'
$ws.Range("AB3").Value = '16:40'

# Row 4
$ws.Range("A4").Value = 111964622
$ws.Range("B4").Value = 89845
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 1209
$ws.Range("F4").Value = 'Rynkskinn'
$ws.Range("G4").Value = 'Phlebia centrifuga'
$ws.Range("H4").Value = 'P.Karst.'
$ws.Range("Q4").Value = 734972.3834676194
$ws.Range("R4").Value = 7088252.533270728
$ws.Range("Z4").Value = '16:12'
$ws.Range("AB4").Value = '16:12'

# Row 5
$ws.Range("A5").Value = 111965883
$ws.Range("B5").Value = 55611
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 102612
$ws.Range("F5").Value = 'Järpe'
$ws.Range("G5").Value = 'Tetrastes bonasia'
$ws.Range("H5").Value = '(Linnaeus, 1758)'
$ws.Range("M5").Value = 'lockläte, övriga läten'
$ws.Range("Q5").Value = 734846.6442297549
$ws.Range("R5").Value = 7088238.22626837
$ws.Range("Z5").Value = '17:05'
$ws.Range("AB5").Value = '17:05'

# Row 6
$ws.Range("A6").Value = 111964847
$ws.Range("B6").Value = 89405
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = 'Ullticka'
$ws.Range("G6").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H6").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("M6").Value = ""
$ws.Range("Q6").Value = 734972.3834676194
$ws.Range("R6").Value = 7088252.533270728
$ws.Range("Z6").Value = '16:12'
$ws.Range("AB6").Value = '16:12'

# Row 7
$ws.Range("A7").Value = 111964632
$ws.Range("B7").Value = 77515
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = 'Garnlav'
$ws.Range("G7").Value = 'Alectoria sarmentosa'
$ws.Range("H7").Value = '(Ach.) Ach.'

# Row 8
$ws.Range("A8").Value = 111964050
$ws.Range("B8").Value = 90065
$ws.Range("D8").Value = 'VU'
$ws.Range("E8").Value = 898
$ws.Range("F8").Value = 'Blackticka'
$ws.Range("G8").Value = 'Steccherinum collabens'
$ws.Range("H8").Value = '(Fr.) Vesterholt'
$ws.Range("Q8").Value = 734893.3330648565
$ws.Range("R8").Value = 7088354.646951701
$ws.Range("Z8").Value = '15:42'
$ws.Range("AB8").Value = '15:42'

# Row 9
$ws.Range("A9").Value = 111965370
$ws.Range("B9").Value = 81248
$ws.Range("E9").Value = 1312
$ws.Range("F9").Value = 'Gammelgransskål'
$ws.Range("G9").Value = 'Pseudographis pinicola'
$ws.Range("H9").Value = '(Nyl.) Rehm'
$ws.Range("Q9").Value = 734939.7547518623
$ws.Range("R9").Value = 7088232.371273324
$ws.Range("Z9").Value = '16:38'
$ws.Range("AB9").Value = '16:38'

# Row 10
$ws.Range("A10").Value = 111964175
$ws.Range("B10").Value = 89423
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 5432
$ws.Range("F10").Value = 'Granticka'
$ws.Range("G10").Value = 'Porodaedalea chrysoloma'
$ws.Range("H10").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q10").Value = 734896.4627943118
$ws.Range("R10").Value = 7088342.483217424
$ws.Range("Z10").Value = '15:42'
$ws.Range("AB10").Value = '15:42'

# Row 11
$ws.Range("A11").Value = 111964863
$ws.Range("B11").Value = 89745
$ws.Range("E11").Value = 2062
$ws.Range("F11").Value = 'Ulltickeporing'
$ws.Range("G11").Value = 'Skeletocutis brevispora'
$ws.Range("H11").Value = 'Niemelä'
$ws.Range("Q11").Value = 734972.3834676194
$ws.Range("R11").Value = 7088252.533270728
$ws.Range("Z11").Value = '16:12'
$ws.Range("AB11").Value = '16:12'
